$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: header labels (species names), no special number format ---
$ws.Range("B23").Value = "malaria"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "toxoplasmosis"
$ws.Range("D23").Value = "chlamidia"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "brucei"
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").Value = "cruzi"
$ws.Range("G23").Value = "leishmania"

# --- Row 24: nucleotide A frequencies ---
$ws.Range("A24").Value = "A"
$ws.Range("B24").Value = 0.40314943767226402
$ws.Range("C24").Value = 0.238634742293792
$ws.Range("D24").Value = 0.29421142444406301
$ws.Range("E24").Value = 0.26903184221300802
$ws.Range("F24").Value = 0.248454666990474
$ws.Range("G24").Value = 0.199678180684498
$ws.Range("B24:G24").NumberFormat = "0.00%"

# --- Row 25: nucleotide C frequencies ---
$ws.Range("A25").Value = "C"
$ws.Range("B25").Value = 0.096649796239937705
$ws.Range("C25").Value = 0.26155470646036
$ws.Range("D25").Value = 0.206453791249848
$ws.Range("E25").Value = 0.232030672793332
$ws.Range("F25").Value = 0.255657303028617
$ws.Range("G25").Value = 0.29988849213646801
$ws.Range("B25:G25").NumberFormat = "0.00%"

# --- Row 26: nucleotide G frequencies ---
$ws.Range("A26").Value = "G"
$ws.Range("B26").Value = 0.096968580375922001
$ws.Range("C26").Value = 0.26122664518829097
$ws.Range("D26").Value = 0.20661877625251901
$ws.Range("E26").Value = 0.230319371291408
$ws.Range("F26").Value = 0.249876237483249
$ws.Range("G26").Value = 0.29730883033559302
$ws.Range("B26:G26").NumberFormat = "0.00%"

# --- Row 27: nucleotide T frequencies ---
$ws.Range("A27").Value = "T"
$ws.Range("B27").Value = 0.40323218571187702
$ws.Range("C27").Value = 0.238583906057558
$ws.Range("D27").Value = 0.29271600805357001
$ws.Range("E27").Value = 0.26861811370225203
$ws.Range("F27").Value = 0.24601179249766
$ws.Range("G27").Value = 0.203124496843441
$ws.Range("B27:G27").NumberFormat = "0.00%"

$ws.Range("R31").Select()

Write-Output "data written"
